$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2613.875
$ws.Cells.Item(98, 9).Value = 2136.077
$ws.Cells.Item(98, 10).Value = 4684.3335
$ws.Cells.Item(98, 11).Value = 2136.077
$ws.Cells.Item(98, 12).Value = 4684.3335
$ws.Cells.Item(98, 13).Value = -638.0770000000002
$ws.Cells.Item(98, 14).Value = -7680.3335
$ws.Cells.Item(101, 8).Value = 375
$ws.Cells.Item(101, 9).Value = 155
$ws.Cells.Item(101, 10).Value = 595
$ws.Cells.Item(101, 11).Value = 465
$ws.Cells.Item(101, 12).Value = 1785
$ws.Cells.Item(101, 13).Value = 1157
$ws.Cells.Item(101, 14).Value = -5029
$ws.Cells.Item(112, 8).Value = 7259.7896
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 7259.7896
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 21779.3688
$ws.Cells.Item(112, 13).Value = ""
$ws.Cells.Item(112, 14).Value = -23995.3688
$ws.Cells.Item(122, 8).Value = 2613.875
$ws.Cells.Item(122, 9).Value = 2136.077
$ws.Cells.Item(122, 10).Value = 4684.3335
$ws.Cells.Item(122, 11).Value = 6408.231000000001
$ws.Cells.Item(122, 12).Value = 14053.0005
$ws.Cells.Item(122, 13).Value = -3958.231000000001
$ws.Cells.Item(122, 14).Value = -18953.0005
$ws.Cells.Item(137, 8).Value = 47178.855
$ws.Cells.Item(137, 9).Value = 104084.664
$ws.Cells.Item(137, 10).Value = 4499.5
$ws.Cells.Item(137, 11).Value = 312253.992
$ws.Cells.Item(137, 12).Value = 13498.5
$ws.Cells.Item(137, 13).Value = -309703.992
$ws.Cells.Item(137, 14).Value = -18598.5
$ws.Cells.Item(138, 8).Value = 2649.976
$ws.Cells.Item(138, 9).Value = 1635.7826
$ws.Cells.Item(138, 10).Value = 3877.6843
$ws.Cells.Item(138, 11).Value = 4907.3478
$ws.Cells.Item(138, 12).Value = 11633.0529
$ws.Cells.Item(138, 13).Value = 232.6522000000004
$ws.Cells.Item(138, 14).Value = -21913.0529

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4849.4688
$ws.Cells.Item(61, 9).Value = 1106.7778
$ws.Cells.Item(61, 10).Value = 25060
$ws.Cells.Item(61, 11).Value = 1106.7778
$ws.Cells.Item(61, 12).Value = 25060
$ws.Cells.Item(61, 13).Value = -894.7778000000001
$ws.Cells.Item(61, 14).Value = -25484
$ws.Cells.Item(97, 8).Value = 1199.44
$ws.Cells.Item(97, 9).Value = 808.2727
$ws.Cells.Item(97, 10).Value = 1506.7858
$ws.Cells.Item(97, 11).Value = 808.2727
$ws.Cells.Item(97, 12).Value = 1506.7858
$ws.Cells.Item(97, 13).Value = -312.2727
$ws.Cells.Item(97, 14).Value = -2498.7858
$ws.Cells.Item(132, 8).Value = 1395.3704
$ws.Cells.Item(132, 9).Value = 1153.125
$ws.Cells.Item(132, 10).Value = 3333.3333
$ws.Cells.Item(132, 11).Value = 3459.375
$ws.Cells.Item(132, 12).Value = 9999.999899999999
$ws.Cells.Item(132, 13).Value = -929.375
$ws.Cells.Item(132, 14).Value = -15059.9999
$ws.Cells.Item(136, 8).Value = 4849.4688
$ws.Cells.Item(136, 9).Value = 1106.7778
$ws.Cells.Item(136, 10).Value = 25060
$ws.Cells.Item(136, 11).Value = 3320.3334
$ws.Cells.Item(136, 12).Value = 75180
$ws.Cells.Item(136, 13).Value = -770.3334000000004
$ws.Cells.Item(136, 14).Value = -80280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1500.6666
$ws.Cells.Item(94, 9).Value = 549.55554
$ws.Cells.Item(94, 10).Value = 2451.7778
$ws.Cells.Item(94, 11).Value = 549.55554
$ws.Cells.Item(94, 12).Value = 2451.7778
$ws.Cells.Item(94, 13).Value = -98.55553999999995
$ws.Cells.Item(94, 14).Value = -3353.7778
$ws.Cells.Item(99, 8).Value = 2442.6667
$ws.Cells.Item(99, 9).Value = 2212.1428
$ws.Cells.Item(99, 10).Value = 3249.5
$ws.Cells.Item(99, 11).Value = 2212.1428
$ws.Cells.Item(99, 12).Value = 3249.5
$ws.Cells.Item(99, 13).Value = -714.1428000000001
$ws.Cells.Item(99, 14).Value = -6245.5
$ws.Cells.Item(105, 8).Value = 11670.889
$ws.Cells.Item(105, 9).Value = 12317.25
$ws.Cells.Item(105, 10).Value = 6500
$ws.Cells.Item(105, 11).Value = 12317.25
$ws.Cells.Item(105, 12).Value = 6500
$ws.Cells.Item(105, 13).Value = -10570.25
$ws.Cells.Item(105, 14).Value = -9994
$ws.Cells.Item(134, 8).Value = 4232.4165
$ws.Cells.Item(134, 9).Value = 4179.2
$ws.Cells.Item(134, 10).Value = 4498.5
$ws.Cells.Item(134, 11).Value = 12537.6
$ws.Cells.Item(134, 12).Value = 13495.5
$ws.Cells.Item(134, 13).Value = -10002.6
$ws.Cells.Item(134, 14).Value = -18565.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6669522.5
$ws.Cells.Item(31, 9).Value = 9092986
$ws.Cells.Item(31, 10).Value = 4999.5
$ws.Cells.Item(31, 11).Value = 9092986
$ws.Cells.Item(31, 12).Value = 4999.5
$ws.Cells.Item(31, 13).Value = -9092691
$ws.Cells.Item(31, 14).Value = -5589.5
$ws.Cells.Item(34, 8).Value = 6669522.5
$ws.Cells.Item(34, 9).Value = 9092986
$ws.Cells.Item(34, 10).Value = 4999.5
$ws.Cells.Item(34, 11).Value = 9092986
$ws.Cells.Item(34, 12).Value = 4999.5
$ws.Cells.Item(34, 13).Value = -9092784
$ws.Cells.Item(34, 14).Value = -5403.5
$ws.Cells.Item(58, 8).Value = 1459.4706
$ws.Cells.Item(58, 9).Value = 1574.0667
$ws.Cells.Item(58, 10).Value = 600
$ws.Cells.Item(58, 11).Value = 1574.0667
$ws.Cells.Item(58, 12).Value = 600
$ws.Cells.Item(58, 13).Value = -1371.0667
$ws.Cells.Item(58, 14).Value = -1006
$ws.Cells.Item(94, 8).Value = 821.8889
$ws.Cells.Item(94, 9).Value = 952.1667
$ws.Cells.Item(94, 10).Value = 561.3333
$ws.Cells.Item(94, 11).Value = 952.1667
$ws.Cells.Item(94, 12).Value = 561.3333
$ws.Cells.Item(94, 13).Value = -501.1667
$ws.Cells.Item(94, 14).Value = -1463.3333
$ws.Cells.Item(107, 8).Value = 1968.4
$ws.Cells.Item(107, 9).Value = 303.625
$ws.Cells.Item(107, 10).Value = 3871
$ws.Cells.Item(107, 11).Value = 303.625
$ws.Cells.Item(107, 12).Value = 3871
$ws.Cells.Item(107, 13).Value = 1616.375
$ws.Cells.Item(107, 14).Value = -7711
$ws.Cells.Item(132, 8).Value = 27786.436
$ws.Cells.Item(132, 9).Value = 31416.75
$ws.Cells.Item(132, 10).Value = 3584.3333
$ws.Cells.Item(132, 11).Value = 94250.25
$ws.Cells.Item(132, 12).Value = 10752.9999
$ws.Cells.Item(132, 13).Value = -91720.25
$ws.Cells.Item(132, 14).Value = -15812.9999
$ws.Cells.Item(134, 8).Value = 2862.7144
$ws.Cells.Item(134, 9).Value = 2528.5
$ws.Cells.Item(134, 10).Value = 3308.3333
$ws.Cells.Item(134, 11).Value = 7585.5
$ws.Cells.Item(134, 12).Value = 9924.999899999999
$ws.Cells.Item(134, 13).Value = -5050.5
$ws.Cells.Item(134, 14).Value = -14994.9999
$ws.Cells.Item(136, 8).Value = 1459.4706
$ws.Cells.Item(136, 9).Value = 1574.0667
$ws.Cells.Item(136, 10).Value = 600
$ws.Cells.Item(136, 11).Value = 4722.2001
$ws.Cells.Item(136, 12).Value = 1800
$ws.Cells.Item(136, 13).Value = -2172.2001
$ws.Cells.Item(136, 14).Value = -6900

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 1832.1333
$ws.Cells.Item(107, 9).Value = 3624.75
$ws.Cells.Item(107, 10).Value = 1180.2727
$ws.Cells.Item(107, 11).Value = 10874.25
$ws.Cells.Item(107, 12).Value = 3540.8181
$ws.Cells.Item(107, 13).Value = -8954.25
$ws.Cells.Item(107, 14).Value = -7380.8181
$ws.Cells.Item(122, 8).Value = 1180.7632
$ws.Cells.Item(122, 9).Value = 950
$ws.Cells.Item(122, 10).Value = 1215.7273
$ws.Cells.Item(122, 11).Value = 8550
$ws.Cells.Item(122, 12).Value = 10941.5457
$ws.Cells.Item(122, 13).Value = -6100
$ws.Cells.Item(122, 14).Value = -15841.5457

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 913.1429000000001
$ws.Cells.Item(22, 9).Value = 499.83334
$ws.Cells.Item(22, 10).Value = 1223.125
$ws.Cells.Item(22, 11).Value = 499.83334
$ws.Cells.Item(22, 12).Value = 1223.125
$ws.Cells.Item(22, 13).Value = -204.83334
$ws.Cells.Item(22, 14).Value = -1813.125
$ws.Cells.Item(27, 8).Value = 913.1429000000001
$ws.Cells.Item(27, 9).Value = 499.83334
$ws.Cells.Item(27, 10).Value = 1223.125
$ws.Cells.Item(27, 11).Value = 499.83334
$ws.Cells.Item(27, 12).Value = 1223.125
$ws.Cells.Item(27, 13).Value = -392.83334
$ws.Cells.Item(27, 14).Value = -1437.125
$ws.Cells.Item(61, 8).Value = 1079.45
$ws.Cells.Item(61, 9).Value = 990.5833
$ws.Cells.Item(61, 10).Value = 1212.75
$ws.Cells.Item(61, 11).Value = 990.5833
$ws.Cells.Item(61, 12).Value = 1212.75
$ws.Cells.Item(61, 13).Value = -788.5833
$ws.Cells.Item(61, 14).Value = -1616.75
$ws.Cells.Item(113, 8).Value = 1079.45
$ws.Cells.Item(113, 9).Value = 990.5833
$ws.Cells.Item(113, 10).Value = 1212.75
$ws.Cells.Item(113, 11).Value = 990.5833
$ws.Cells.Item(113, 12).Value = 1212.75
$ws.Cells.Item(113, 13).Value = 1179.4167
$ws.Cells.Item(113, 14).Value = -5552.75
$ws.Cells.Item(122, 8).Value = 3734.3215
$ws.Cells.Item(122, 9).Value = 3632.9565
$ws.Cells.Item(122, 10).Value = 4200.6
$ws.Cells.Item(122, 11).Value = 10898.8695
$ws.Cells.Item(122, 12).Value = 12601.8
$ws.Cells.Item(122, 13).Value = -8448.869499999999
$ws.Cells.Item(122, 14).Value = -17501.8
$ws.Cells.Item(132, 8).Value = 2931.8918
$ws.Cells.Item(132, 9).Value = 2999.639
$ws.Cells.Item(132, 10).Value = 493
$ws.Cells.Item(132, 11).Value = 8998.917000000001
$ws.Cells.Item(132, 12).Value = 1479
$ws.Cells.Item(132, 13).Value = -6468.917000000001
$ws.Cells.Item(132, 14).Value = -6539
$ws.Cells.Item(136, 8).Value = 5192.6
$ws.Cells.Item(136, 9).Value = 5172.3335
$ws.Cells.Item(136, 10).Value = 5375
$ws.Cells.Item(136, 11).Value = 15517.0005
$ws.Cells.Item(136, 12).Value = 16125
$ws.Cells.Item(136, 13).Value = -12967.0005
$ws.Cells.Item(136, 14).Value = -21225

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 32950.27
$ws.Cells.Item(132, 9).Value = 33868.28
$ws.Cells.Item(132, 10).Value = 10000
$ws.Cells.Item(132, 11).Value = 101604.84
$ws.Cells.Item(132, 12).Value = 30000
$ws.Cells.Item(132, 13).Value = -99074.84
$ws.Cells.Item(132, 14).Value = -35060
$ws.Cells.Item(136, 8).Value = 28139.846
$ws.Cells.Item(136, 9).Value = 30005.666
$ws.Cells.Item(136, 10).Value = 5750
$ws.Cells.Item(136, 11).Value = 90016.99800000001
$ws.Cells.Item(136, 12).Value = 17250
$ws.Cells.Item(136, 13).Value = -87466.99800000001
$ws.Cells.Item(136, 14).Value = -22350
